# Code updation on orderconfirmation and receipt voucher test page
#
# Updates the "month / year / day" style trio of fields that were tested
# with Jan/2023/20 (or /4) to Feb/2023/28 on the OrderDetailPage sheet, and
# bumps the matching "day" field from 1 to 18 on the OrderConfirmationPage
# and StoreHeadConfirmationPage sheets (to line up with the already-updated
# CreateCustomerOrderNo reference sheet). Also moves the "active" sheet /
# selection bookmarks that Excel persists when the workbook is saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# OrderDetailPage: chequeMonth/Year/Day-style triples -> Feb / 2023 / 28
# ---------------------------------------------------------------------
$wsOrderDetail = $wb.Worksheets.Item("OrderDetailPage")

$wsOrderDetail.Range("H2").Value = "Feb"
$wsOrderDetail.Range("J2").Value = "28"

$wsOrderDetail.Range("W2").Value = "Feb"
$wsOrderDetail.Range("Y2").Value = "28"

$wsOrderDetail.Range("AU2").Value = "Feb"
$wsOrderDetail.Range("AW2").Value = "28"

$wsOrderDetail.Activate() | Out-Null
$wsOrderDetail.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------
# OrderConfirmationPage: day field 1 -> 18 (twice), widen the new column
# ---------------------------------------------------------------------
$wsOrderConfirmation = $wb.Worksheets.Item("OrderConfirmationPage")

$wsOrderConfirmation.Range("E2").Value = "18"
$wsOrderConfirmation.Range("H2").Value = "18"
$wsOrderConfirmation.Columns.Item(10).ColumnWidth = 13.109375

$wsOrderConfirmation.Activate() | Out-Null
$wsOrderConfirmation.Range("H2").Select() | Out-Null

# ---------------------------------------------------------------------
# StoreHeadConfirmationPage: day field 1 -> 18 (twice); this becomes the
# workbook's active sheet/tab on save
# ---------------------------------------------------------------------
$wsStoreHead = $wb.Worksheets.Item("StoreHeadConfirmationPage")

$wsStoreHead.Range("E2").Value = "18"
$wsStoreHead.Range("H2").Value = "18"

$wsStoreHead.Activate() | Out-Null
$wsStoreHead.Range("H7").Select() | Out-Null

# ---------------------------------------------------------------------
# CreateCustomerOrderNo was the previously-active sheet; it keeps its own
# selection (A2) but is no longer the active tab now that
# StoreHeadConfirmationPage has been activated above.
# ---------------------------------------------------------------------
